# Updated translation patch to support the huge v2 update.
# Column B becomes the translation column (mirroring column A by default),
# with the previously scattered C/D translations ("Lily", "Shina", "Lime")
# consolidated into B for rows 2-4. Columns C and D are no longer used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 46

# Default: column B mirrors column A's text for every data row.
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value2
}

# Rows 2-4 carry real translations (previously stored in columns C/D).
$ws.Range("B2").Value = "Lily"
$ws.Range("B3").Value = "Shina"
$ws.Range("B4").Value = "Lime"

# Columns C and D are no longer needed.
$ws.Range("C1:D46").ClearContents()

# Re-fit row heights so multi-line cell text doesn't leave a stale
# explicit/custom row height behind.
$ws.Rows("1:" + $lastRow).AutoFit()
